$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1805.5714
$ws.Range("I6").Value = 1827.4
$ws.Range("K6").Value = 5482.200000000001
$ws.Range("M6").Value = -5370.200000000001
$ws.Range("H9").Value = 1332676.5
$ws.Range("I9").Value = 322.625
$ws.Range("K9").Value = 322.625
$ws.Range("M9").Value = -153.625
$ws.Range("H132").Value = 1838.4445
$ws.Range("I132").Value = 1505.75
$ws.Range("K132").Value = 4517.25
$ws.Range("M132").Value = -1987.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3202.8
$ws.Range("I45").Value = 2500
$ws.Range("J45").Value = 3671.3333
$ws.Range("K45").Value = 2500
$ws.Range("L45").Value = 3671.3333
$ws.Range("M45").Value = -2123
$ws.Range("N45").Value = -4425.3333
$ws.Range("H61").Value = 9740.207
$ws.Range("I61").Value = 7663.1055
$ws.Range("K61").Value = 7663.1055
$ws.Range("M61").Value = -7451.1055
$ws.Range("H74").Value = 3368.3333
$ws.Range("I74").Value = 1112
$ws.Range("K74").Value = 1112
$ws.Range("M74").Value = -238
$ws.Range("H77").Value = 3368.3333
$ws.Range("I77").Value = 1112
$ws.Range("K77").Value = 5560
$ws.Range("M77").Value = -1192
$ws.Range("H92").Value = 39999.8
$ws.Range("J92").Value = 39999.8
$ws.Range("L92").Value = 39999.8
$ws.Range("N92").Value = -44991.8
$ws.Range("H102").Value = 2935.5
$ws.Range("I102").Value = 2935.5
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2935.5
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -1313.5
$ws.Range("H104").Value = 4517.3335
$ws.Range("J104").Value = 4380
$ws.Range("L104").Value = 4380
$ws.Range("N104").Value = -11368
$ws.Range("H122").Value = 3237.25
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3237.25
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").Value = 9711.75
$ws.Range("N122").Value = -14611.75
$ws.Range("H132").Value = 6388.394
$ws.Range("I132").Value = 4575.9565
$ws.Range("J132").Value = 10557
$ws.Range("K132").Value = 13727.8695
$ws.Range("L132").Value = 31671
$ws.Range("M132").Value = -11197.8695
$ws.Range("N132").Value = -36731
$ws.Range("H136").Value = 9740.207
$ws.Range("I136").Value = 7663.1055
$ws.Range("K136").Value = 22989.3165
$ws.Range("M136").Value = -20439.3165

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 96566.28999999999
$ws.Range("I86").Value = 1251.0834
$ws.Range("J86").Value = 223653.22
$ws.Range("K86").Value = 1251.0834
$ws.Range("L86").Value = 223653.22
$ws.Range("M86").Value = -128.0834
$ws.Range("N86").Value = -225899.22
$ws.Range("H89").Value = 96566.28999999999
$ws.Range("I89").Value = 1251.0834
$ws.Range("J89").Value = 223653.22
$ws.Range("K89").Value = 6255.416999999999
$ws.Range("L89").Value = 1118266.1
$ws.Range("M89").Value = -639.4169999999995
$ws.Range("N89").Value = -1129498.1
$ws.Range("H94").Value = 1228.7097
$ws.Range("I94").Value = 1013.0455
$ws.Range("K94").Value = 1013.0455
$ws.Range("M94").Value = -562.0454999999999
$ws.Range("H105").Value = 250007500
$ws.Range("I105").Value = 250007500
$ws.Range("K105").Value = 250007500
$ws.Range("M105").Value = -250005753
$ws.Range("H134").Value = 3384.7058
$ws.Range("I134").Value = 1792.6216
$ws.Range("K134").Value = 5377.864799999999
$ws.Range("M134").Value = -2842.864799999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4832.5654
$ws.Range("I58").Value = 3959.8
$ws.Range("K58").Value = 3959.8
$ws.Range("M58").Value = -3756.8
$ws.Range("H134").Value = 5603.7617
$ws.Range("I134").Value = 4837.8237
$ws.Range("K134").Value = 14513.4711
$ws.Range("M134").Value = -11978.4711
$ws.Range("H136").Value = 4832.5654
$ws.Range("I136").Value = 3959.8
$ws.Range("K136").Value = 11879.4
$ws.Range("M136").Value = -9329.400000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 80.8421
$ws.Range("I12").Value = 4.75
$ws.Range("J12").Value = 101.13333
$ws.Range("K12").Value = 14.25
$ws.Range("L12").Value = 303.39999
$ws.Range("M12").Value = 158.75
$ws.Range("N12").Value = -649.39999
$ws.Range("H34").Value = 663.3125
$ws.Range("J34").Value = 10003
$ws.Range("L34").Value = 30009
$ws.Range("N34").Value = -30177
$ws.Range("H99").Value = 2303.5715
$ws.Range("I99").Value = 562.5
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 1687.5
$ws.Range("L99").Value = 9000
$ws.Range("M99").Value = 558.5
$ws.Range("N99").Value = -13492

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 1500
$ws.Range("I9").Value = 1500
$ws.Range("K9").Value = 1500
$ws.Range("M9").Value = -1330
$ws.Range("H70").Value = 4666.5
$ws.Range("I70").Value = 4666.5
$ws.Range("K70").Value = 4666.5
$ws.Range("M70").Value = -4396.5
$ws.Range("H73").Value = 4666.5
$ws.Range("I73").Value = 4666.5
$ws.Range("K73").Value = 4666.5
$ws.Range("M73").Value = -3730.5
$ws.Range("H95").Value = 34998
$ws.Range("J95").Value = 34998
$ws.Range("L95").Value = 34998
$ws.Range("N95").Value = -40490
$ws.Range("H97").Value = 1386.15
$ws.Range("I97").Value = 1451.6428
$ws.Range("K97").Value = 1451.6428
$ws.Range("M97").Value = -955.6428000000001
$ws.Range("H102").Value = 994.1111
$ws.Range("I102").Value = 855.1429000000001
$ws.Range("J102").Value = 1480.5
$ws.Range("K102").Value = 855.1429000000001
$ws.Range("L102").Value = 1480.5
$ws.Range("M102").Value = 766.8570999999999
$ws.Range("N102").Value = -4724.5
$ws.Range("H122").Value = 1371.8235
$ws.Range("I122").Value = 1532
$ws.Range("K122").Value = 4596
$ws.Range("M122").Value = -2146
$ws.Range("H132").Value = 4533.1724
$ws.Range("I132").Value = 2293.7273
$ws.Range("K132").Value = 6881.1819
$ws.Range("M132").Value = -4351.1819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2615.111
$ws.Range("I7").Value = 2027
$ws.Range("K7").Value = 2027
$ws.Range("M7").Value = -1915
$ws.Range("H22").Value = 3027.0852
$ws.Range("I22").Value = 2054.9473
$ws.Range("J22").Value = 3686.75
$ws.Range("K22").Value = 2054.9473
$ws.Range("L22").Value = 3686.75
$ws.Range("M22").Value = -1759.9473
$ws.Range("N22").Value = -4276.75
$ws.Range("H27").Value = 3027.0852
$ws.Range("I27").Value = 2054.9473
$ws.Range("J27").Value = 3686.75
$ws.Range("K27").Value = 2054.9473
$ws.Range("L27").Value = 3686.75
$ws.Range("M27").Value = -1947.9473
$ws.Range("N27").Value = -3900.75
$ws.Range("H126").Value = 2615.111
$ws.Range("I126").Value = 2027
$ws.Range("K126").Value = 6081
$ws.Range("M126").Value = -3611
$ws.Range("H132").Value = 5645.2856
$ws.Range("I132").Value = 3996
$ws.Range("J132").Value = 7294.5713
$ws.Range("K132").Value = 11988
$ws.Range("L132").Value = 21883.7139
$ws.Range("M132").Value = -9458
$ws.Range("N132").Value = -26943.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 69047.5
$ws.Range("I70").Value = 60095
$ws.Range("J70").Value = 78000
$ws.Range("K70").Value = 60095
$ws.Range("L70").Value = 78000
$ws.Range("M70").Value = -59780
$ws.Range("N70").Value = -78630
$ws.Range("H73").Value = 69047.5
$ws.Range("I73").Value = 60095
$ws.Range("J73").Value = 78000
$ws.Range("K73").Value = 60095
$ws.Range("L73").Value = 78000
$ws.Range("M73").Value = -59003
$ws.Range("N73").Value = -80184
$ws.Range("H107").Value = 871.3182
$ws.Range("I107").Value = 531
$ws.Range("K107").Value = 1593
$ws.Range("M107").Value = 327
$ws.Range("H126").Value = 4652.9414
$ws.Range("I126").Value = 4166.385
$ws.Range("K126").Value = 12499.155
$ws.Range("M126").Value = -10029.155
$ws.Range("H132").Value = 5067.72
$ws.Range("I132").Value = 3502.8572
$ws.Range("J132").Value = 7059.364
$ws.Range("K132").Value = 10508.5716
$ws.Range("L132").Value = 21178.092
$ws.Range("M132").Value = -7978.571599999999
$ws.Range("N132").Value = -26238.092
